# Updates the "Price" (column D) and "Volume(1h)" (column E) values for the
# crypto rows on Sheet1, matching the refreshed data from the scheduled
# GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D="38.421.29"; E="  +1.86%  "},
    @{Row=3; D="2.078.98"; E="  +2.15%  "},
    @{Row=4; D=$null; E="  -0.04%  "},
    @{Row=5; D="228.71"; E="  +0.58%  "},
    @{Row=6; D=$null; E="  +0.42%  "},
    @{Row=7; D=$null; E="  +0.33%  "},
    @{Row=8; D=$null; E="  +0.02%  "},
    @{Row=9; D="0.381"; E="  +1.42%  "},
    @{Row=10; D="0.0832"; E="  +0.74%  "},
    @{Row=11; D=$null; E="  -0.62%  "},
    @{Row=12; D="2.387.79"; E="  +2.34%  "},
    @{Row=13; D="14.85"; E="  +2.17%  "},
    @{Row=14; D="22.40"; E="  +6.13%  "},
    @{Row=15; D=$null; E="  +0.91%  "},
    @{Row=16; D=$null; E="  +2.20%  "},
    @{Row=17; D="2.081.78"; E="  +2.78%  "},
    @{Row=18; D="38.371.05"; E="  +1.90%  "},
    @{Row=19; D="71.15"; E="  +2.52%  "},
    @{Row=20; D="6.01"; E="  +0.86%  "},
    @{Row=21; D=$null; E="  +1.33%  "},
    @{Row=22; D="224.98"; E="  +0.46%  "},
    @{Row=23; D=$null; E="  -0.11%  "},
    @{Row=24; D=$null; E="  -0.13%  "},
    @{Row=25; D="2.34"; E="  +2.66%  "},
    @{Row=26; D="169.88"; E="  +1.27%  "},
    @{Row=27; D="9.39"; E="  +0.41%  "},
    @{Row=28; D=$null; E="  +6.54%  "},
    @{Row=29; D=$null; E="  +1.42%  "},
    @{Row=30; D=$null; E="  +8.72%  "},
    @{Row=31; D=$null; E="  -0.64%  "},
    @{Row=32; D="2.32"; E="  +4.91%  "},
    @{Row=33; D="4.77"; E="  +6.73%  "},
    @{Row=34; D=$null; E="  +3.01%  "},
    @{Row=35; D="0.0606"; E="  +0.18%  "},
    @{Row=36; D=$null; E="  +1.05%  "},
    @{Row=37; D="6.34"; E="  -2.78%  "},
    @{Row=38; D="3.53"; E="  +4.02%  "},
    @{Row=39; D=$null; E="  +0.02%  "},
    @{Row=40; D="18.30"; E="  +2.20%  "},
    @{Row=41; D="1.540.95"; E="  +0.89%  "},
    @{Row=42; D="100.11"; E="  +3.17%  "},
    @{Row=43; D=$null; E="  +2.20%  "},
    @{Row=44; D=$null; E="  +1.40%  "},
    @{Row=45; D=$null; E="  -1.05%  "},
    @{Row=46; D="7.69"; E="  +9.10%  "},
    @{Row=47; D=$null; E="  +0.32%  "},
    @{Row=48; D=$null; E="  -1.69%  "},
    @{Row=49; D=$null; E="  +2.03%  "},
    @{Row=50; D="2.98"; E="  +1.64%  "},
    @{Row=51; D="2.277.33"; E="  +2.47%  "}
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Column D values (e.g. "6.01", "228.71") look numeric to Excel's
        # smart-typing, which would otherwise silently coerce them into
        # doubles (losing trailing zeros / introducing float noise). Force
        # the cell to Text, assign the literal string, then clear the
        # number-format override so the cell's style stays untouched
        # (matching the original workbook, which has no explicit style on
        # these data cells).
        $cellD = $ws.Cells.Item($u.Row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
        $cellD.ClearFormats()
    }

    if ($null -ne $u.E) {
        $cellE = $ws.Cells.Item($u.Row, 5)
        $cellE.NumberFormat = "@"
        $cellE.Value = $u.E
        $cellE.ClearFormats()
    }
}
